# Generate Report for Handback
#
# The 162628d1-... file's handback transform failed because the handback
# archive's file name didn't match the original handoff file name. Update
# the localization-status report: flip the status from "Ready for handoff"
# to "Handback transform failed" (Overview + both locale tabs) and record
# the explanatory message in the "Error Detail" column for each locale.

$wb = $excel.ActiveWorkbook

$errorDetailZhCn = "Handback file name: aw53idec.5cz is different with handoff file name: 162628d1-f4d5-41d0-9b12-3d89575ba153.b708c197d5749d57ecf1a835888fc8b65738f115.zh-cn."
$errorDetailDeDe = "Handback file name: aw53idec.5cz is different with handoff file name: 162628d1-f4d5-41d0-9b12-3d89575ba153.b708c197d5749d57ecf1a835888fc8b65738f115.de-de."

# Overview tab: zh-cn / de-de status columns for the 162628d1-... row.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# zh-cn tab: Status + Error Detail for the 162628d1-... row, and widen the
# Error Detail column (column P, the 16th) to 40 characters so the new
# message is readable. (39.17 is what lands on a stored width of 40 once
# Excel quantizes ColumnWidth to whole pixels on save.)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = $errorDetailZhCn
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# de-de tab: Status + Error Detail for the 162628d1-... row, and widen the
# Error Detail column the same way.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = $errorDetailDeDe
$dede.Columns.Item(16).ColumnWidth = 39.17
